$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Caribbean_frogs_distribution_da")

# Determine the extent of data in column C (Island names)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -eq "Jamaica") {
        $cell.Value = "Hispaniola"
    } elseif ($val -eq "Hispaniola") {
        $cell.Value = "Jamaica"
    }
}

# Update the active selection to match the recorded edit location
$ws.Range("E18").Select()
